$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "1.001") are not
# auto-converted to numbers by Excel's input parser; cleared afterwards so no
# residual style attribute is left on the cells.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = '28.292.42'
$ws.Range("E2").Value = '  +2.16%  '
$ws.Range("D3").Value = '1.815.05'
$ws.Range("E3").Value = '  +3.43%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '325.53'
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '0.4374'
$ws.Range("E7").Value = '  +1.62%  '
$ws.Range("E8").Value = '  +0.66%  '
$ws.Range("D9").Value = '44.74'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("D10").Value = '0.07665'
$ws.Range("E10").Value = '  +2.40%  '
$ws.Range("D11").Value = '1.140'
$ws.Range("E11").Value = '  +1.83%  '
$ws.Range("D12").Value = '0.9996'
$ws.Range("E12").Value = '  -0.08%  '
$ws.Range("D13").Value = '21.96'
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").Value = '6.299'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").Value = '7.491'
$ws.Range("E15").Value = '  +3.34%  '
$ws.Range("D16").Value = '1.828.61'
$ws.Range("E16").Value = '  +4.62%  '
$ws.Range("D17").Value = '95.05'
$ws.Range("E17").Value = '  +8.22%  '
$ws.Range("D18").Value = '0.00001079'
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("D19").Value = '0.06496'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '17.39'
$ws.Range("E21").Value = '  +1.60%  '
$ws.Range("D22").Value = '6.241'
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("D23").Value = '28.304.72'
$ws.Range("E23").Value = '  +2.10%  '
$ws.Range("D24").Value = '11.55'
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").Value = '2.110'
$ws.Range("E25").Value = '  -9.31%  '
$ws.Range("D26").Value = '161.32'
$ws.Range("E26").Value = '  +5.66%  '
$ws.Range("D27").Value = '20.73'
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("D28").Value = '2.032.35'
$ws.Range("E28").Value = '  +4.33%  '
$ws.Range("D29").Value = '2.275'
$ws.Range("E29").Value = '  -3.78%  '
$ws.Range("D30").Value = '129.20'
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("D31").Value = '1.207'
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").Value = '6.009'
$ws.Range("E32").Value = '  +5.01%  '
$ws.Range("D33").Value = '0.09141'
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("D34").Value = '3.573'
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("D35").Value = '12.95'
$ws.Range("E35").Value = '  +1.93%  '
$ws.Range("D36").Value = '0.02363'
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").Value = '5.223'
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("D38").Value = '0.2173'
$ws.Range("E38").Value = '  +0.96%  '
$ws.Range("D39").Value = '0.6587'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("D41").Value = '1.190'
$ws.Range("E41").Value = '  -0.62%  '
$ws.Range("D42").Value = '8.076'
$ws.Range("E42").Value = '  +1.68%  '
$ws.Range("D43").Value = '1.429'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '0.9984'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").Value = '13.78'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = '0.6100'
$ws.Range("E46").Value = '  +2.65%  '
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").Value = '125.40'
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("E49").Value = '  +2.07%  '
$ws.Range("D50").Value = '1.157'
$ws.Range("E50").Value = '  +2.99%  '
$ws.Range("D51").Value = '0.06994'
$ws.Range("E51").Value = '  +1.33%  '

$colD.ClearFormats()
